$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Copy()
$ws.Range("A33").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = $false

$ws.Range("A33").Value = 46003
$ws.Range("B33").Value = 64

$ws.Range("A33:B33").Select()
